$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D: old D (value 4) shifts to E, C is untouched.
$ws.Columns("D").Insert()

# Clear out the old B2 value (SUBID count of 1).
$ws.Range("B2").ClearContents()

# Fill in the newly-inserted column D and a further new column F on row 2.
$ws.Range("D2").Value = 3
$ws.Range("F2").Value = 6

# Tweak the column widths slightly (rank/epoch-binning related column resize).
$ws.Columns("A").ColumnWidth = 5.59
$ws.Columns("B").ColumnWidth = 11.59

# Move the active selection (visualization tweak - last touched cell).
$ws.Range("L15").Select() | Out-Null
